# Add a new "Permission Name" column (F) to the Item_Code sheet, mirroring the
# Menu/Button Name information so each permission row also carries a
# human-readable name (menu name for the "Menu" rows, action name otherwise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell - copy the formatting of the adjacent header cell (bold)
# and then set its text.
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F2").Value = "Permission Name"

# New column width to match the other descriptive columns (closest
# achievable value to the source width of 29.77734375 character-units,
# given the engine's internal rounding of ColumnWidth).
$ws.Columns("F").ColumnWidth = 29

# Data rows: "Menu" header rows get the Menu/Group name, action rows mirror
# the Button Name column.
$ws.Range("F4").Value = "Dashboard"

$ws.Range("F6").Value = "Company"
$ws.Range("F7").Value = "Add"
$ws.Range("F8").Value = "Edit"
$ws.Range("F9").Value = "Delete"
$ws.Range("F10").Value = "Export"

$ws.Range("F12").Value = "Permission Group"
$ws.Range("F13").Value = "Add"
$ws.Range("F14").Value = "Edit"
$ws.Range("F15").Value = "Delete"
$ws.Range("F16").Value = "Export"

$ws.Range("F18").Value = "User"
$ws.Range("F19").Value = "Add"
$ws.Range("F20").Value = "Edit"
$ws.Range("F21").Value = "Delete"
$ws.Range("F22").Value = "Export"

# Match the final selection left by the author.
$ws.Range("F22").Select()
